$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Insert a brand-new worksheet "2022-Q4" right after "总计" (first sheet)
# ---------------------------------------------------------------------------
$zongji = $wb.Worksheets.Item("总计")

$newSheet = $wb.Worksheets.Add($null, $zongji)
$newSheet.Name = "2022-Q4"

# Re-fetch the "2022-Q3" sheet by name now that sheet order/indices changed
# because of the Add() above (stale index-based refs silently fail to carry
# formats through PasteSpecial).
$q3 = $wb.Worksheets.Item("2022-Q3")

# Match page margins used by every other quarter sheet
$newSheet.PageSetup.LeftMargin   = 0.75 * 72
$newSheet.PageSetup.RightMargin  = 0.75 * 72
$newSheet.PageSetup.TopMargin    = 1    * 72
$newSheet.PageSetup.BottomMargin = 1    * 72
$newSheet.PageSetup.HeaderMargin = 0.5  * 72
$newSheet.PageSetup.FooterMargin = 0.5  * 72

# Headers (B1:H1) - copy the bold/bordered/centered style used on every
# other quarter sheet's header row, then fill in the text.
$q3.Range("B1:H1").Copy()
$newSheet.Range("B1:H1").PasteSpecial(-4122)   # xlPasteFormats

$headers = @("基金代码","基金名称","基金规模","股票总仓位","仓位占比","持有市值(亿元)","仓位排名")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $newSheet.Cells.Item(1, 2 + $i).Value = $headers[$i]
}

# Data rows
$data = @(
    @("001672","国寿安保智慧生活股票","10.45","90.56","2.58","0.2696",9),
    @("004818","国寿安保目标策略灵活配置混合A","2.76","59.92","3.07","0.0847",5),
    @("004819","国寿安保目标策略灵活配置混合C","1.30","59.92","3.07","0.0399",5),
    @("090019","大成景恒混合A","1.18","93.72","1.58","0.0186",9),
    @("006038","大成景恒混合C","0.89","93.72","1.58","0.0141",9)
)

# A-column style (bold/bordered/centered, same as header) used for the
# little numeric row-index column on every quarter sheet.
$q3.Range("A2").Copy()
$newSheet.Range("A2:A6").PasteSpecial(-4122)

for ($r = 0; $r -lt $data.Length; $r++) {
    $row = 2 + $r
    $rowData = $data[$r]

    $newSheet.Cells.Item($row, 1).Value = $r          # A: 0-based index (numeric, already styled)

    # B,D,E,F,G look like numbers but must stay TEXT (inline string) cells,
    # exactly like on the other quarter sheets - force text format, enter
    # the value, then drop back to the default "Normal" style so no stray
    # number-format survives on the cell.
    foreach ($col in 2,4,5,6,7) {
        $cell = $newSheet.Cells.Item($row, $col)
        $cell.NumberFormat = "@"
        $cell.Value = $rowData[$col - 2]
        $cell.Style = "Normal"
    }

    # C: plain textual fund name, stores fine as text without any trick.
    $newSheet.Cells.Item($row, 3).Value = $rowData[1]

    # H: numeric rank.
    $newSheet.Cells.Item($row, 8).Value = $rowData[6]
}

# ---------------------------------------------------------------------------
# 2) Update "总计": insert a new row 2 for "2022-Q4", pushing every other
#    quarter down by one row, and renumber the A-column index.
# ---------------------------------------------------------------------------
$zongji = $wb.Worksheets.Item("总计")

$zongji.Rows(2).Insert()
$zongji.Range("B2:D2").Style = "Normal"

# Give A2 the same numeric/bordered style as the rest of the A column.
$zongji.Range("A3").Copy()
$zongji.Range("A2").PasteSpecial(-4122)

$zongji.Range("A2").Value = 0
$zongji.Range("B2").Value = "2022-Q4"
$zongji.Range("C2").Value = 5
$zongji.Range("D2").Value = 0.43

for ($r = 3; $r -le 7; $r++) {
    $zongji.Cells.Item($r, 1).Value = $r - 2
}

Write-Host "2022-Q4 sheet added and 总计 updated"
